$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (col D) and Volume(1h) (col E) figures.
# Rows 11-15, 20-21 and 37-39 were also re-ranked by the scraper, so the
# coin name (col B) and link (col C) cells move together with their new data.
#
# Price cells are plain text in the source data (e.g. thousand-separated
# "27.080.35" or values with significant trailing zeros like "1.010"), so we
# force text formatting before assigning them to stop Excel from silently
# re-interpreting them as numbers, then restore the default "Normal" style so
# no extra formatting is left behind on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '27.080.35'
$ws.Cells.Item(2, 5).Value = '  +0.10%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.828.73'
$ws.Cells.Item(3, 5).Value = '  +0.29%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '1.011'
$ws.Cells.Item(4, 5).Value = '  +0.48%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '312.37'
$ws.Cells.Item(5, 5).Value = '  +0.05%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '1.009'
$ws.Cells.Item(6, 5).Value = '  +0.32%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.4635'
$ws.Cells.Item(7, 5).Value = '  -1.13%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.3706'
$ws.Cells.Item(8, 5).Value = '  +1.58%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.07357'
$ws.Cells.Item(9, 5).Value = '  -0.37%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.8728'
$ws.Cells.Item(10, 5).Value = '  -0.55%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(11, 4) '1.914.87'
$ws.Cells.Item(11, 5).Value = '  +1.71%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Cells.Item(12, 4) '19.89'
$ws.Cells.Item(12, 5).Value = '  -1.59%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Cells.Item(13, 4) '0.07838'
$ws.Cells.Item(13, 5).Value = '  +3.48%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(14, 4) '5.345'
$ws.Cells.Item(14, 5).Value = '  -0.38%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Chainlink'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Cells.Item(15, 4) '6.570'
$ws.Cells.Item(15, 5).Value = '  +0.76%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '92.01'
$ws.Cells.Item(16, 5).Value = '  -0.85%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '1.011'
$ws.Cells.Item(17, 5).Value = '  +0.58%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '0.000008871'
$ws.Cells.Item(18, 5).Value = '  +1.78%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.16%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Avalanche'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Cells.Item(20, 4) '14.68'
$ws.Cells.Item(20, 5).Value = '  +0.66%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'WrappedBTC'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Cells.Item(21, 4) '27.038.49'
$ws.Cells.Item(21, 5).Value = '  -1.57%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '5.145'
$ws.Cells.Item(22, 5).Value = '  -1.61%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '10.64'
$ws.Cells.Item(23, 5).Value = '  +0.33%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '2.140.06'
$ws.Cells.Item(24, 5).Value = '  +2.94%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '152.50'
$ws.Cells.Item(25, 5).Value = '  +0.75%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '1.849'
$ws.Cells.Item(26, 5).Value = '  -1.83%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '18.47'
$ws.Cells.Item(27, 5).Value = '  +0.10%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '2.084'
$ws.Cells.Item(28, 5).Value = '  -1.97%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '5.086'
$ws.Cells.Item(29, 5).Value = '  -1.53%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '115.70'
$ws.Cells.Item(30, 5).Value = '  -0.51%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '0.08874'
$ws.Cells.Item(31, 5).Value = '  -0.44%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '2.974'
$ws.Cells.Item(32, 5).Value = '  +1.16%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '0.7319'
$ws.Cells.Item(33, 5).Value = '  -1.53%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '4.448'
$ws.Cells.Item(34, 5).Value = '  -1.33%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '1.136'
$ws.Cells.Item(35, 5).Value = '  -2.03%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '2.485'
$ws.Cells.Item(36, 5).Value = '  -8.22%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(37, 4) '0.05267'
$ws.Cells.Item(37, 5).Value = '  -0.36%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Cells.Item(38, 4) '1.073'
$ws.Cells.Item(38, 5).Value = '  -1.23%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(39, 4) '0.01947'
$ws.Cells.Item(39, 5).Value = '  +0.99%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '2.928'
$ws.Cells.Item(40, 5).Value = '  -0.03%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '7.167'
$ws.Cells.Item(41, 5).Value = '  -1.71%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '0.5193'
$ws.Cells.Item(42, 5).Value = '  -1.09%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '0.8873'
$ws.Cells.Item(43, 5).Value = '  -11.80%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.1632'
$ws.Cells.Item(44, 5).Value = '  -0.50%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '8.226'
$ws.Cells.Item(45, 5).Value = '  -1.60%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '0.4842'
$ws.Cells.Item(46, 5).Value = '  -1.19%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '1.010'
$ws.Cells.Item(47, 5).Value = '  +0.36%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '10.18'
$ws.Cells.Item(48, 5).Value = '  -1.51%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '102.37'
$ws.Cells.Item(49, 5).Value = '  -1.90%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '1.627'
$ws.Cells.Item(50, 5).Value = '  -1.35%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '0.06226'
$ws.Cells.Item(51, 5).Value = '  -0.59%  '
